# Update the header cell C1 to fix the "Change in Status" label:
#   old: "Change in Status (2007 to 2017"
#   new: "Change in status (2007 to 2017)"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C1").Value = "Change in status (2007 to 2017)"

# Update the active selection to match the saved state (B2 selected).
$ws.Range("B2").Select()
